$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "69.289.59"
$ws.Range("E2").Value = "  -1.62%  "
$ws.Range("D3").Value = "3.491.78"
$ws.Range("E3").Value = "  -3.20%  "
$ws.Range("E4").Value = "  -0.03%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "580.25"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -0.33%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "182.26"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -4.12%  "
$ws.Range("D7").Value = "3.480.73"
$ws.Range("E7").Value = "  -3.45%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.609"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  -3.59%  "
$ws.Range("E9").Value = "  +0.03%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.195"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +5.94%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.640"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -3.31%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "53.69"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  -4.41%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.0000302"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  -3.66%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "9.36"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -3.52%  "
$ws.Range("D15").Value = "4.032.74"
$ws.Range("E15").Value = "  -3.80%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "19.18"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  -3.31%  "
$ws.Range("D17").Value = "69.195.52"
$ws.Range("E17").Value = "  -1.70%  "
$ws.Range("D18").Value = "3.500.15"
$ws.Range("E18").Value = "  -2.99%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "12.23"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -3.27%  "
$ws.Range("E20").Value = "  -1.52%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "531.38"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +7.97%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "1.00"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -4.21%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "18.44"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -4.69%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "4.50"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +3.10%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "4.84"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -1.30%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "95.38"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -1.51%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "11.00"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -0.40%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "2.94"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -1.92%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "9.04"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -3.66%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "31.81"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -1.47%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "7.21"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -4.56%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "12.34"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +0.64%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "63.61"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -3.41%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.112"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -4.77%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "537.08"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -7.68%  "
$ws.Range("B36").Value = "TheGraph"
$ws.Range("C36").Value = "https://coinranking.com/coin/qhd1biQ7M+thegraph-grt"
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.405"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +0.99%  "
$ws.Range("B37").Value = "Fetch.AI"
$ws.Range("C37").Value = "https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet"
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "3.08"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +3.61%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "37.84"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -3.14%  "
$ws.Range("B39").Value = "Dai"
$ws.Range("C39").Value = "https://coinranking.com/coin/MoTuySvg7+dai-dai"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "1.00"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -0.10%  "
$ws.Range("D40").Value = "0.0₃0749"
$ws.Range("E40").Value = "  -8.28%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.134"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -2.39%  "
$ws.Range("D42").Value = "3.337.15"
$ws.Range("E42").Value = "  +3.37%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "3.35"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -4.17%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "3.52"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +3.78%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "2.98"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -8.78%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "2.94"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -4.22%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.0433"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -2.92%  "
$ws.Range("E48").Value = "  -3.85%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "8.94"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -8.56%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.997"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -0.22%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "136.82"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +0.29%  "

Write-Host "Applied all cell updates"